# Hand in documents for the final
# Update hours-spent figures in the "Members" table and the surrounding
# narrative paragraphs of the final report.

$d = $word.ActiveDocument

# --- "Members" table: "Hours spent on project" column ------------------
$table = $d.Tables.Item(1)

# Fries, Lea: 120 -> 133
$table.Cell(2, 4).Range.Find.Execute("120", $true, $false, $false, $false, $false, $true, 1, $false, "133", 2)

# Karcher, Lukas: 179 -> 180
$table.Cell(3, 4).Range.Find.Execute("179", $true, $false, $false, $false, $false, $true, 1, $false, "180", 2)

# Kertzscher, Tim: 176 -> 177
$table.Cell(4, 4).Range.Find.Execute("176", $true, $false, $false, $false, $false, $true, 1, $false, "177", 2)

# Rickel, Jan: 194 -> 200
$table.Cell(5, 4).Range.Find.Execute("194", $true, $false, $false, $false, $false, $true, 1, $false, "200", 2)

# --- Narrative paragraphs below the table --------------------------------

# "Out of these 669 hours, 301 were spent ..." -> "... 690 hours, 309 were spent ..."
$d.Content.Find.Execute("Out of these 669 hours, 301 were", $true, $false, $false, $false, $false, $true, 1, $false, "Out of these 690 hours, 309 were", 2)

# Date the hours were accumulated: 06/16/2019 -> 06/17/2019
$d.Content.Find.Execute("06/16/2019", $true, $false, $false, $false, $false, $true, 1, $false, "06/17/2019", 2)

# Outstanding hours: "Another 10h" -> "Another 2h"
$d.Content.Find.Execute("Another 10h", $true, $false, $false, $false, $false, $true, 1, $false, "Another 2h", 2)
